# Inception_Phase.pptx edit
#
# Summary of the change (derived from the OOXML diff):
#   The "Initial use case" slide (old slide 49, containing a centered
#   title placeholder + a screenshot picture) is split into two slides:
#     - a new slide 49: just the (resized/repositioned) title "Initial use case"
#       acting as a section/title card, no picture.
#     - slide 50 (the old slide 49, pushed down by one): keeps the picture
#       (slightly repositioned) and the title (repositioned to match the
#       new title card), plus a new "Fig. 1: Initial use case" caption
#       text box added underneath the picture.
#   Every following slide shifts down by one position (no content changes).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Duplicate the "Initial use case" slide (old index 49) and move the
#    duplicate in front of the original -> this becomes the new slide 49
#    (title card) while the original becomes slide 50 (full content).
# ---------------------------------------------------------------------
$srcIndex = 49
$original = $p.Slides.Item($srcIndex)

$dupRange = $original.Duplicate()
$titleCard = $dupRange.Item(1)
$titleCard.MoveTo($srcIndex)

# After the move:
#   slide 49 -> $titleCard  (new, title-only)
#   slide 50 -> $original   (old content, now needs pic reposition + caption)

# ---------------------------------------------------------------------
# 2. New slide 49 ("title card"): remove the picture, resize/reposition
#    the title text.
# ---------------------------------------------------------------------
$titleCardPic = $titleCard.Shapes.Item("Picture 3")
$titleCardPic.Delete()

$titleCardTitle = $titleCard.Shapes.Item("Title 1")
$titleCardTitle.Left = 2573835 / 12700
$titleCardTitle.Top = 613896 / 12700
$titleCardTitle.Width = 6574544 / 12700
$titleCardTitle.Height = 365335 / 12700
$titleCardTitle.TextFrame.AutoSize = 0
$titleCardTitle.TextFrame.TextRange.Text = "Initial use case"
$titleCardTitle.TextFrame.TextRange.Font.Size = 32
$titleCardTitle.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# 3. Slide 50 (the old "Initial use case" content slide): reposition the
#    title to match the title-card layout, reposition the picture, and
#    add the new "Fig. 1: Initial use case" caption text box.
# ---------------------------------------------------------------------
$content = $original

$contentTitle = $content.Shapes.Item("Title 1")
$contentTitle.Left = 2573835 / 12700
$contentTitle.Top = 613896 / 12700
$contentTitle.Width = 6574544 / 12700
$contentTitle.Height = 365335 / 12700
$contentTitle.TextFrame.AutoSize = 0
$contentTitle.TextFrame.TextRange.Text = "Initial use case"
$contentTitle.TextFrame.TextRange.Font.Size = 32
$contentTitle.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$contentPic = $content.Shapes.Item("Picture 3")
$contentPic.Left = 2175173 / 12700
$contentPic.Top = 1190252 / 12700
$contentPic.Width = 7371869 / 12700
$contentPic.Height = 5143840 / 12700

$caption = $content.Shapes.AddTextbox(1, 1887474 / 12700, 6294438 / 12700, 6574544 / 12700, 365335 / 12700)
$caption.Name = "Title 1"
$caption.TextFrame.WordWrap = $true
$caption.TextFrame.AutoSize = 2
$caption.TextFrame.TextRange.Text = "Fig. 1: Initial use case"
$caption.TextFrame.TextRange.Font.Size = 16
$caption.TextFrame.TextRange.ParagraphFormat.Alignment = 2

Write-Output ("Final slide count: " + $p.Slides.Count)
